$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1771.6666
$ws.Range("J6").Value = 40
$ws.Range("L6").Value = 120
$ws.Range("N6").Value = -344
$ws.Range("H9").Value = 767
$ws.Range("I9").Value = 767
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 767
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -598
$ws.Range("N9").Value = ""
$ws.Range("H12").Value = 221.66667
$ws.Range("I12").Value = 80
$ws.Range("J12").Value = 363.33334
$ws.Range("K12").Value = 80
$ws.Range("L12").Value = 363.33334
$ws.Range("M12").Value = 90
$ws.Range("N12").Value = -703.33334
$ws.Range("H21").Value = 29999.857
$ws.Range("I21").Value = 29999
$ws.Range("K21").Value = 29999
$ws.Range("M21").Value = -29531
$ws.Range("H23").Value = 29999.857
$ws.Range("I23").Value = 29999
$ws.Range("K23").Value = 29999
$ws.Range("M23").Value = -29765
$ws.Range("H29").Value = 1411.5
$ws.Range("I29").Value = 436.8
$ws.Range("J29").Value = 1854.5454
$ws.Range("K29").Value = 1310.4
$ws.Range("L29").Value = 5563.6362
$ws.Range("M29").Value = -1029.4
$ws.Range("N29").Value = -6125.6362
$ws.Range("H38").Value = 293.33334
$ws.Range("I38").Value = 152
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 456
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = -84
$ws.Range("N38").Value = -3744
$ws.Range("H58").Value = 3814.2856
$ws.Range("I58").Value = 233.33333
$ws.Range("J58").Value = 6500
$ws.Range("K58").Value = 699.99999
$ws.Range("L58").Value = 19500
$ws.Range("M58").Value = -549.99999
$ws.Range("N58").Value = -19800
$ws.Range("H129").Value = 1684498.2
$ws.Range("I129").Value = 451.44446
$ws.Range("J129").Value = 2850376.8
$ws.Range("K129").Value = 1354.33338
$ws.Range("L129").Value = 8551130.399999999
$ws.Range("M129").Value = 3645.66662
$ws.Range("N129").Value = -8561130.399999999
$ws.Range("H137").Value = 1073.2333
$ws.Range("I137").Value = 778.4783
$ws.Range("J137").Value = 2041.7142
$ws.Range("K137").Value = 2335.4349
$ws.Range("L137").Value = 6125.142599999999
$ws.Range("M137").Value = 214.5650999999998
$ws.Range("N137").Value = -11225.1426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 2000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = -2230
$ws.Range("H25").Value = 2866.6667
$ws.Range("I25").Value = 800
$ws.Range("J25").Value = 3900
$ws.Range("K25").Value = 800
$ws.Range("L25").Value = 3900
$ws.Range("M25").Value = -398
$ws.Range("N25").Value = -4704
$ws.Range("H28").Value = 20354.2
$ws.Range("I28").Value = 6323.6665
$ws.Range("K28").Value = 6323.6665
$ws.Range("M28").Value = -6131.6665
$ws.Range("H41").Value = 19200
$ws.Range("I41").Value = 15000
$ws.Range("J41").Value = 20600
$ws.Range("K41").Value = 15000
$ws.Range("L41").Value = 20600
$ws.Range("M41").Value = -14586
$ws.Range("N41").Value = -21428
$ws.Range("H61").Value = 932.2
$ws.Range("I61").Value = 816.1177
$ws.Range("J61").Value = 1590
$ws.Range("K61").Value = 816.1177
$ws.Range("L61").Value = 1590
$ws.Range("M61").Value = -604.1177
$ws.Range("N61").Value = -2014
$ws.Range("H99").Value = 20354.2
$ws.Range("I99").Value = 6323.6665
$ws.Range("K99").Value = 6323.6665
$ws.Range("M99").Value = -3328.6665
$ws.Range("H107").Value = 30000
$ws.Range("J107").Value = 30000
$ws.Range("L107").Value = 30000
$ws.Range("N107").Value = -37680
$ws.Range("H109").Value = 29800
$ws.Range("J109").Value = 29800
$ws.Range("L109").Value = 29800
$ws.Range("N109").Value = -32574
$ws.Range("H136").Value = 932.2
$ws.Range("I136").Value = 816.1177
$ws.Range("J136").Value = 1590
$ws.Range("K136").Value = 2448.3531
$ws.Range("L136").Value = 4770
$ws.Range("M136").Value = 101.6468999999997
$ws.Range("N136").Value = -9870

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = ""
$ws.Range("H11").Value = 46935
$ws.Range("I11").Value = 795
$ws.Range("J11").Value = 70005
$ws.Range("K11").Value = 795
$ws.Range("L11").Value = 70005
$ws.Range("M11").Value = -655
$ws.Range("N11").Value = -70285
$ws.Range("H37").Value = 3369.1428
$ws.Range("I37").Value = 2175.3333
$ws.Range("J37").Value = 4264.5
$ws.Range("K37").Value = 2175.3333
$ws.Range("L37").Value = 4264.5
$ws.Range("M37").Value = -2038.3333
$ws.Range("N37").Value = -4538.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2091.975
$ws.Range("I31").Value = 2126.4595
$ws.Range("K31").Value = 2126.4595
$ws.Range("M31").Value = -1831.4595
$ws.Range("H34").Value = 2091.975
$ws.Range("I34").Value = 2126.4595
$ws.Range("K34").Value = 2126.4595
$ws.Range("M34").Value = -1924.4595

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2069.077
$ws.Range("I4").Value = 733.3333
$ws.Range("J4").Value = 2243.3044
$ws.Range("K4").Value = 2199.9999
$ws.Range("L4").Value = 6729.9132
$ws.Range("M4").Value = -2087.9999
$ws.Range("N4").Value = -6953.9132
$ws.Range("H34").Value = 1503.8
$ws.Range("I34").Value = 539.6667
$ws.Range("J34").Value = 2950
$ws.Range("K34").Value = 1619.0001
$ws.Range("L34").Value = 8850
$ws.Range("M34").Value = -1535.0001
$ws.Range("N34").Value = -9018
$ws.Range("H39").Value = 2436.3157
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 2436.3157
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 7308.9471
$ws.Range("M39").Value = ""
$ws.Range("N39").Value = -7896.9471
$ws.Range("H55").Value = 46489.91
$ws.Range("I55").Value = 333800
$ws.Range("J55").Value = 1125.1578
$ws.Range("K55").Value = 1001400
$ws.Range("L55").Value = 3375.4734
$ws.Range("M55").Value = -1001223
$ws.Range("N55").Value = -3729.4734
$ws.Range("H87").Value = 13800
$ws.Range("I87").Value = 5950
$ws.Range("J87").Value = 21650
$ws.Range("K87").Value = 17850
$ws.Range("L87").Value = 64950
$ws.Range("M87").Value = -16602
$ws.Range("N87").Value = -67446
$ws.Range("H90").Value = 13800
$ws.Range("I90").Value = 5950
$ws.Range("J90").Value = 21650
$ws.Range("K90").Value = 53550
$ws.Range("L90").Value = 194850
$ws.Range("M90").Value = -47310
$ws.Range("N90").Value = -207330
$ws.Range("H121").Value = 1250.6666
$ws.Range("I121").Value = 1397.5
$ws.Range("J121").Value = 1133.2
$ws.Range("K121").Value = 4192.5
$ws.Range("L121").Value = 3399.6
$ws.Range("M121").Value = -2882.5
$ws.Range("N121").Value = -6019.6
$ws.Range("H129").Value = 51668.75
$ws.Range("I129").Value = 1598.5714
$ws.Range("J129").Value = 78629.62
$ws.Range("K129").Value = 4795.7142
$ws.Range("L129").Value = 235888.86
$ws.Range("M129").Value = 204.2857999999997
$ws.Range("N129").Value = -245888.86
$ws.Range("H131").Value = 790.26
$ws.Range("I131").Value = 418
$ws.Range("J131").Value = 809.85266
$ws.Range("K131").Value = 1254
$ws.Range("L131").Value = 2429.55798
$ws.Range("M131").Value = 3786
$ws.Range("N131").Value = -12509.55798
$ws.Range("H139").Value = 2020.8182
$ws.Range("I139").Value = 2027.6666
$ws.Range("J139").Value = 1990
$ws.Range("K139").Value = 6082.9998
$ws.Range("L139").Value = 5970
$ws.Range("M139").Value = -942.9997999999996
$ws.Range("N139").Value = -16250
$ws.Range("H141").Value = 6586.45
$ws.Range("J141").Value = 9031.25
$ws.Range("L141").Value = 27093.75
$ws.Range("N141").Value = -37453.75
